$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 506; this shifts the existing 506.. block down by 3
# and (as observed) inherits formatting (e.g. the date style on column D) from the
# row immediately below the insertion point.
$ws.Rows("506:508").Insert()

# Common (unchanged-across-the-block) column values, copied from the surrounding rows.
$mercadoId = 8
$mercado = "Terminal La Palmera de La Serena"
$region = "Coquimbo"
$codreg = 4
$tipo = "Fruta"
$productoId = 100101
$producto = "Berries"
$categoriaId = 100112025
$categoria = "Frutilla"
$variedad = "Sin especificar"
$unidad = "`$/bandeja 7 kilos"
$origen = "Provincia de Melipilla"
$kgUnidad = 7

function Set-FrutillaRow {
    param($row, $fecha, $calidad, $volumen, $precioMin, $precioMax, $precioProm, $precioKg)

    $ws.Cells.Item($row, 1).Value = $mercadoId
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $fecha
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $tipo
    $ws.Cells.Item($row, 7).Value = $productoId
    $ws.Cells.Item($row, 8).Value = $producto
    $ws.Cells.Item($row, 9).Value = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $calidad
    $ws.Cells.Item($row, 13).Value = $volumen
    $ws.Cells.Item($row, 14).Value = $precioMin
    $ws.Cells.Item($row, 15).Value = $precioMax
    $ws.Cells.Item($row, 16).Value = $precioProm
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $precioKg
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}

# New week of Frutilla prices (Terminal La Palmera de La Serena), dated 44578.
Set-FrutillaRow 506 44578 "Especial" 400 11500 12000 11750 1679
Set-FrutillaRow 507 44578 "Primera"  400 9500  10000 9750  1393
Set-FrutillaRow 508 44578 "Segunda"  360 7500  8000  7750  1107
